$wb = $excel.ActiveWorkbook

# The Custid value (column L, row 2) is updated on the three sheets that
# share this same test customer id (NewCust, DeleteCust, EditCust).
$sheetNames = @("NewCust", "DeleteCust", "EditCust")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("L2").Value = "30876"
}
